$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the ASV_ID / Species_name / Common_name / Category values between
# row 42 and row 43 (columns A-D). Columns E and F remain unchanged.
$a42 = $ws.Range("A42").Value2
$b42 = $ws.Range("B42").Value2
$c42 = $ws.Range("C42").Value2
$d42 = $ws.Range("D42").Value2

$a43 = $ws.Range("A43").Value2
$b43 = $ws.Range("B43").Value2
$c43 = $ws.Range("C43").Value2
$d43 = $ws.Range("D43").Value2

$ws.Range("A42").Value = $a43
$ws.Range("B42").Value = $b43
$ws.Range("C42").Value = $c43
$ws.Range("D42").Value = $d43

$ws.Range("A43").Value = $a42
$ws.Range("B43").Value = $b42
$ws.Range("C43").Value = $c42
$ws.Range("D43").Value = $d42
